# Insert a new data row at row 165 (this shifts rows 165..242 down to 166..243,
# growing the sheet's dimension from A1:R242 to A1:R243).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(165).Insert()

# The newly inserted row 165 is blank. Populate it by duplicating the row that is
# now directly below it (row 166, which holds what used to be row 165's data), then
# overwrite just the cells that actually differ for the new record.
$ws.Range("A166:R166").Copy()
$ws.Range("A165:R165").PasteSpecial()

# Apply the new record's own values on top of the duplicated row.
$ws.Cells.Item(165, 4).Value = 44553   # D165 - Fecha
$ws.Cells.Item(165, 10).Value = 55     # J165 - Volumen
$ws.Cells.Item(165, 11).Value = 5000   # K165 - Precio minimo
$ws.Cells.Item(165, 12).Value = 5000   # L165 - Precio maximo
$ws.Cells.Item(165, 13).Value = 5000   # M165 - Precio promedio ponderado
$ws.Cells.Item(165, 16).Value = 1667   # P165 - Precio $/Kg
